$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.2484426666666667
$ws.Range("H2").Value = 0.745328
$ws.Range("I2").Value = 0.1396403772415532
$ws.Range("J2").Value = 0.1396403772415532
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.171693666666666
$ws.Range("N2").Value = 12.515081
$ws.Range("O2").Value = 0.1077921033402881
$ws.Range("P2").Value = 0.1077921033402881
$ws.Range("Q2").Value = 1.036426699063111
$ws.Range("R2").Value = 9.327840291567998
$ws.Range("S2").Value = 0.01505212997409832
$ws.Range("T2").Value = 0.01505212997409833

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.2484426666666667
$ws.Range("H3").Value = 0.745328
$ws.Range("I3").Value = 0.1396403772415532
$ws.Range("J3").Value = 0.1396403772415532
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 20.39394566666667
$ws.Range("N3").Value = 61.181837
$ws.Range("O3").Value = 0.5269577477327285
$ws.Range("P3").Value = 0.5269577477327286
$ws.Range("Q3").Value = 5.066726245281778
$ws.Range("R3").Value = 45.600536207536
$ws.Range("S3").Value = 0.07358457868375745
$ws.Range("T3").Value = 0.07358457868375747

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.2484426666666667
$ws.Range("H4").Value = 0.745328
$ws.Range("I4").Value = 0.1396403772415532
$ws.Range("J4").Value = 0.1396403772415532
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.13565266666667
$ws.Range("N4").Value = 42.406958
$ws.Range("O4").Value = 0.3652501489269833
$ws.Range("P4").Value = 0.3652501489269833
$ws.Range("Q4").Value = 3.511899243580445
$ws.Range("R4").Value = 31.607093192224
$ws.Range("S4").Value = 0.05100366858369744
$ws.Range("T4").Value = 0.05100366858369745

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.055305333333333
$ws.Range("H5").Value = 3.165916
$ws.Range("I5").Value = 0.5931478551122046
$ws.Range("J5").Value = 0.5931478551122047
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.171693666666666
$ws.Range("N5").Value = 12.515081
$ws.Range("O5").Value = 0.1077921033402881
$ws.Range("P5").Value = 0.1077921033402881
$ws.Range("Q5").Value = 4.402410575466222
$ws.Range("R5").Value = 39.621695179196
$ws.Range("S5").Value = 0.063936654894325
$ws.Range("T5").Value = 0.06393665489432503

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.055305333333333
$ws.Range("H6").Value = 3.165916
$ws.Range("I6").Value = 0.5931478551122046
$ws.Range("J6").Value = 0.5931478551122047
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 20.39394566666667
$ws.Range("N6").Value = 61.181837
$ws.Range("O6").Value = 0.5269577477327285
$ws.Range("P6").Value = 0.5269577477327286
$ws.Range("Q6").Value = 21.52183962974356
$ws.Range("R6").Value = 193.696556667692
$ws.Range("S6").Value = 0.3125638578024261
$ws.Range("T6").Value = 0.3125638578024262

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.055305333333333
$ws.Range("H7").Value = 3.165916
$ws.Range("I7").Value = 0.5931478551122046
$ws.Range("J7").Value = 0.5931478551122047
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.13565266666667
$ws.Range("N7").Value = 42.406958
$ws.Range("O7").Value = 0.3652501489269833
$ws.Range("P7").Value = 0.3652501489269833
$ws.Range("Q7").Value = 14.91742964928089
$ws.Range("R7").Value = 134.256866843528
$ws.Range("S7").Value = 0.2166473424154534
$ws.Range("T7").Value = 0.2166473424154535

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.4754126666666667
$ws.Range("H8").Value = 1.426238
$ws.Range("I8").Value = 0.2672117676462422
$ws.Range("J8").Value = 0.2672117676462422
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.171693666666666
$ws.Range("N8").Value = 12.515081
$ws.Range("O8").Value = 0.1077921033402881
$ws.Range("P8").Value = 0.1077921033402881
$ws.Range("Q8").Value = 1.983276010586444
$ws.Range("R8").Value = 17.849484095278
$ws.Range("S8").Value = 0.0288033184718648
$ws.Range("T8").Value = 0.0288033184718648

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.4754126666666667
$ws.Range("H9").Value = 1.426238
$ws.Range("I9").Value = 0.2672117676462422
$ws.Range("J9").Value = 0.2672117676462422
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 20.39394566666667
$ws.Range("N9").Value = 61.181837
$ws.Range("O9").Value = 0.5269577477327285
$ws.Range("P9").Value = 0.5269577477327286
$ws.Range("Q9").Value = 9.695540093245112
$ws.Range("R9").Value = 87.25986083920601
$ws.Range("S9").Value = 0.140809311246545
$ws.Range("T9").Value = 0.140809311246545

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4754126666666667
$ws.Range("H10").Value = 1.426238
$ws.Range("I10").Value = 0.2672117676462422
$ws.Range("J10").Value = 0.2672117676462422
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.13565266666667
$ws.Range("N10").Value = 42.406958
$ws.Range("O10").Value = 0.3652501489269833
$ws.Range("P10").Value = 0.3652501489269833
$ws.Range("Q10").Value = 6.720268329333779
$ws.Range("R10").Value = 60.48241496400401
$ws.Range("S10").Value = 0.09759913792783241
$ws.Range("T10").Value = 0.09759913792783244

